$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 4
$ws.Range("G2").Value = 4.6
$ws.Range("H2").Value = 2.02
$ws.Range("I2").Value = 2.16
$ws.Range("P2").Value = 1.81
$ws.Range("Q2").Value = 2.04
$ws.Range("R2").Value = 1.3
$ws.Range("S2").Value = 3.75
$ws.Range("T2").Value = 1.84
$ws.Range("U2").Value = 2
$ws.Range("V2").Value = 1.86
$ws.Range("W2").Value = 1.28
$ws.Range("X2").Value = 15.5
$ws.Range("Y2").Value = 10.5
$ws.Range("Z2").Value = 15.5
$ws.Range("AA2").Value = 32
$ws.Range("AB2").Value = 17.5
$ws.Range("AC2").Value = 9.4
$ws.Range("AE2").Value = 29
$ws.Range("AF2").Value = 38
$ws.Range("AG2").Value = 22
$ws.Range("AH2").Value = 23
$ws.Range("AK2").Value = 70
$ws.Range("AL2").Value = 80
$ws.Range("AM2").Value = 140
$ws.Range("AN2").Value = 80
$ws.Range("AO2").Value = 22
$ws.Range("F3").Value = 1.46
$ws.Range("F4").Value = 1.94
$ws.Range("G4").Value = 1.95
$ws.Range("H4").Value = 4.4
$ws.Range("J4").Value = 3.8
$ws.Range("O4").Value = 1.3
$ws.Range("Q4").Value = 1.89
$ws.Range("R4").Value = 1.41
$ws.Range("S4").Value = 3.25
$ws.Range("T4").Value = 1.8
$ws.Range("U4").Value = 2.16
$ws.Range("AB4").Value = 9.6
$ws.Range("AI4").Value = 60
$ws.Range("AM4").Value = 120
$ws.Range("H5").Value = 12.5
$ws.Range("J5").Value = 5.5
$ws.Range("K5").Value = 5.6
$ws.Range("R5").Value = 1.38
$ws.Range("U5").Value = 1.64
$ws.Range("X5").Value = 16.5
$ws.Range("AB5").Value = 6.8
$ws.Range("AF5").Value = 6.8
$ws.Range("AK5").Value = 16.5
$ws.Range("G6").Value = 3.15
$ws.Range("I6").Value = 3.05
$ws.Range("F7").Value = 1.57
$ws.Range("G7").Value = 1.7
$ws.Range("H7").Value = 5.7
$ws.Range("I7").Value = 7.6
$ws.Range("J7").Value = 4
$ws.Range("K7").Value = 4.8
$ws.Range("N7").Value = 4.1
$ws.Range("O7").Value = 1.25
$ws.Range("P7").Value = 2.1
$ws.Range("Q7").Value = 1.63
$ws.Range("R7").Value = 1.44
$ws.Range("S7").Value = 2.84
$ws.Range("T7").Value = 1.8
$ws.Range("U7").Value = 2
$ws.Range("X7").Value = 23
$ws.Range("Y7").Value = 27
$ws.Range("AC7").Value = 12
$ws.Range("AD7").Value = 29
$ws.Range("AF7").Value = 12.5
$ws.Range("AH7").Value = 26
$ws.Range("AJ7").Value = 19
$ws.Range("AN7").Value = 10
$ws.Range("F8").Value = 11.5
$ws.Range("G8").Value = 12
$ws.Range("H8").Value = 1.27
$ws.Range("I8").Value = 1.28
$ws.Range("J8").Value = 7.4
$ws.Range("P8").Value = 4.1
$ws.Range("S8").Value = 1.72
$ws.Range("T8").Value = 1.63
$ws.Range("U8").Value = 2.52
$ws.Range("X8").Value = 65
$ws.Range("Y8").Value = 20
$ws.Range("Z8").Value = 13.5
$ws.Range("AB8").Value = 75
$ws.Range("AF8").Value = 150
$ws.Range("AG8").Value = 44
$ws.Range("AJ8").Value = 380
$ws.Range("AK8").Value = 150
$ws.Range("AM8").Value = 80
$ws.Range("AN8").Value = 85
$ws.Range("F9").Value = 1.91
$ws.Range("I9").Value = 4.3
$ws.Range("K9").Value = 4.2
$ws.Range("Q9").Value = 1.64
$ws.Range("X9").Value = 23
$ws.Range("Y9").Value = 21
$ws.Range("AA9").Value = 95
$ws.Range("AG9").Value = 11
$ws.Range("AI9").Value = 55
$ws.Range("AK9").Value = 18
$ws.Range("AM9").Value = 80
$ws.Range("AN9").Value = 9
$ws.Range("AO9").Value = 34
$ws.Range("F10").Value = 1.81
$ws.Range("G10").Value = 1.82
$ws.Range("H10").Value = 4.5
$ws.Range("I10").Value = 4.6
$ws.Range("P10").Value = 2.44
$ws.Range("T10").Value = 1.65
$ws.Range("AF10").Value = 13
$ws.Range("AJ10").Value = 21
$ws.Range("AK10").Value = 17.5
$ws.Range("AN10").Value = 8.6
$ws.Range("U11").Value = 1.52
$ws.Range("Z11").Value = 310
$ws.Range("H12").Value = 8.8
$ws.Range("Q12").Value = 1.83
$ws.Range("F13").Value = 6.6
$ws.Range("J13").Value = 4.6
$ws.Range("N13").Value = 4.9
$ws.Range("O13").Value = 1.19
$ws.Range("P13").Value = 2.38
$ws.Range("Q13").Value = 1.61
$ws.Range("R13").Value = 1.54
$ws.Range("S13").Value = 2.52
$ws.Range("U13").Value = 2.08
$ws.Range("X13").Value = 26
$ws.Range("Y13").Value = 12.5
$ws.Range("AB13").Value = 30
$ws.Range("AC13").Value = 13.5
$ws.Range("AD13").Value = 12
$ws.Range("AE13").Value = 18
$ws.Range("AH13").Value = 25
$ws.Range("AI13").Value = 32
$ws.Range("F14").Value = 1.35
$ws.Range("G14").Value = 1.42
$ws.Range("I14").Value = 12
$ws.Range("J14").Value = 5.1
$ws.Range("K14").Value = 6.2
$ws.Range("M14").Value = 1.03
$ws.Range("N14").Value = 5.1
$ws.Range("P14").Value = 2.36
$ws.Range("Q14").Value = 1.6
$ws.Range("R14").Value = 1.51
$ws.Range("S14").Value = 2.4
$ws.Range("T14").Value = 1.91
$ws.Range("U14").Value = 1.91
$ws.Range("V14").Value = 1.09
$ws.Range("W14").Value = 3.35
$ws.Range("X14").Value = 29
$ws.Range("Y14").Value = 42
$ws.Range("AA14").Value = 390
$ws.Range("AB14").Value = 12
$ws.Range("AC14").Value = 15.5
$ws.Range("AD14").Value = 44
$ws.Range("AG14").Value = 13
$ws.Range("AK14").Value = 17.5
$ws.Range("AM14").Value = 160
$ws.Range("AN14").Value = 6.4
